# Auto-generated edit script
# Applies the Marilith_Profits market-data refresh:
# updates currentAveragePrice/-NQ/-HQ, LevePriceNQ/HQ and profit columns
# (H-N) across all 8 sheets to match the latest scheduled-runner pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 999
$ws.Range("I18").Value = 999
$ws.Range("K18").Value = 999
$ws.Range("M18").Value = -715
$ws.Range("H51").Value = 2662.2222
$ws.Range("I51").Value = 1992.5
$ws.Range("J51").Value = 2853.5715
$ws.Range("K51").Value = 1992.5
$ws.Range("L51").Value = 2853.5715
$ws.Range("M51").Value = -1508.5
$ws.Range("N51").Value = -3821.5715
$ws.Range("H62").Value = 5259.8
$ws.Range("I62").Value = 4650
$ws.Range("J62").Value = 5666.3335
$ws.Range("K62").Value = 4650
$ws.Range("L62").Value = 5666.3335
$ws.Range("M62").Value = -4026
$ws.Range("N62").Value = -6914.3335
$ws.Range("H64").Value = 4999.778
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H65").Value = 5259.8
$ws.Range("I65").Value = 4650
$ws.Range("J65").Value = 5666.3335
$ws.Range("K65").Value = 23250
$ws.Range("L65").Value = 28331.6675
$ws.Range("M65").Value = -20130
$ws.Range("N65").Value = -34571.6675
$ws.Range("H67").Value = 4999.778
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H76").Value = 3325.5
$ws.Range("I76").Value = 2613.25
$ws.Range("J76").Value = 4750
$ws.Range("K76").Value = 2613.25
$ws.Range("L76").Value = 4750
$ws.Range("M76").Value = -2298.25
$ws.Range("N76").Value = -5380
$ws.Range("H79").Value = 3325.5
$ws.Range("I79").Value = 2613.25
$ws.Range("J79").Value = 4750
$ws.Range("K79").Value = 2613.25
$ws.Range("L79").Value = 4750
$ws.Range("M79").Value = -1521.25
$ws.Range("N79").Value = -6934
$ws.Range("H99").Value = 1316.2858
$ws.Range("J99").Value = 1567.75
$ws.Range("L99").Value = 4703.25
$ws.Range("N99").Value = -7699.25
$ws.Range("H103").Value = 994.75
$ws.Range("J103").Value = 994.75
$ws.Range("L103").Value = 2984.25
$ws.Range("N103").Value = -4156.25
$ws.Range("H111").Value = 976
$ws.Range("I111").Value = 830
$ws.Range("J111").Value = 1073.3334
$ws.Range("K111").Value = 2490
$ws.Range("L111").Value = 3220.0002
$ws.Range("M111").Value = 577
$ws.Range("N111").Value = -9354.0002
$ws.Range("H112").Value = 3849.5
$ws.Range("H127").Value = 1469.5555
$ws.Range("I127").Value = 778.25
$ws.Range("K127").Value = 2334.75
$ws.Range("M127").Value = 2625.25
$ws.Range("H138").Value = 2604.1428
$ws.Range("I138").Value = 1950.8182
$ws.Range("J138").Value = 4999.6665
$ws.Range("K138").Value = 5852.4546
$ws.Range("L138").Value = 14998.9995
$ws.Range("M138").Value = -712.4546
$ws.Range("N138").Value = -25278.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2167.1667
$ws.Range("I110").Value = 2112.7778
$ws.Range("J110").Value = 2330.3333
$ws.Range("K110").Value = 2112.7778
$ws.Range("L110").Value = 2330.3333
$ws.Range("M110").Value = -67.77779999999984
$ws.Range("N110").Value = -6420.3333
$ws.Range("H132").Value = 1701.5883
$ws.Range("I132").Value = 1701.5883
$ws.Range("K132").Value = 5104.7649
$ws.Range("M132").Value = -2574.7649

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 6020.5
$ws.Range("I54").Value = 6020.5
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 6020.5
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -5536.5
$ws.Range("N54").ClearContents()
$ws.Range("H86").Value = 2179.8572
$ws.Range("I86").Value = 2176.5
$ws.Range("K86").Value = 2176.5
$ws.Range("M86").Value = -1053.5
$ws.Range("H89").Value = 2179.8572
$ws.Range("I89").Value = 2176.5
$ws.Range("K89").Value = 10882.5
$ws.Range("M89").Value = -5266.5
$ws.Range("H107").Value = 1053.6666
$ws.Range("I107").Value = 1064.2
$ws.Range("J107").Value = 1001
$ws.Range("K107").Value = 1064.2
$ws.Range("L107").Value = 1001
$ws.Range("M107").Value = 855.8
$ws.Range("N107").Value = -4841

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 619.8
$ws.Range("I22").Value = 549.6667
$ws.Range("J22").Value = 725
$ws.Range("K22").Value = 549.6667
$ws.Range("L22").Value = 725
$ws.Range("M22").Value = -199.6667
$ws.Range("N22").Value = -1425
$ws.Range("H31").Value = 2299.2666
$ws.Range("I31").Value = 1819.4117
$ws.Range("K31").Value = 1819.4117
$ws.Range("M31").Value = -1524.4117
$ws.Range("H34").Value = 2299.2666
$ws.Range("I34").Value = 1819.4117
$ws.Range("K34").Value = 1819.4117
$ws.Range("M34").Value = -1617.4117
$ws.Range("H105").Value = 1748.25
$ws.Range("I105").Value = 1831
$ws.Range("K105").Value = 1831
$ws.Range("M105").Value = -84
$ws.Range("H107").Value = 1494.125
$ws.Range("I107").Value = 1489.6666
$ws.Range("J107").Value = 1496.8
$ws.Range("K107").Value = 1489.6666
$ws.Range("L107").Value = 1496.8
$ws.Range("M107").Value = 430.3334
$ws.Range("N107").Value = -5336.8
$ws.Range("H134").Value = 1044.3334
$ws.Range("I134").Value = 1057.091
$ws.Range("J134").Value = 904
$ws.Range("K134").Value = 3171.273
$ws.Range("L134").Value = 2712
$ws.Range("M134").Value = -636.2729999999997
$ws.Range("N134").Value = -7782

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 289.75
$ws.Range("I13").Value = 451.8
$ws.Range("K13").Value = 1355.4
$ws.Range("M13").Value = -1187.4
$ws.Range("H68").Value = 1407.2
$ws.Range("J68").Value = 800
$ws.Range("L68").Value = 2400
$ws.Range("N68").Value = -4022
$ws.Range("H71").Value = 1407.2
$ws.Range("J71").Value = 800
$ws.Range("L71").Value = 7200
$ws.Range("N71").Value = -15312
$ws.Range("H80").Value = 14499.75
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 14499.75
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 43499.25
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -45371.25
$ws.Range("H83").Value = 14499.75
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 14499.75
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 130497.75
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -139857.75
$ws.Range("H138").Value = 5842.625
$ws.Range("I138").Value = 5178
$ws.Range("K138").Value = 15534
$ws.Range("M138").Value = -10394

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1337.5
$ws.Range("I97").Value = 1125
$ws.Range("K97").Value = 1125
$ws.Range("M97").Value = -629
$ws.Range("H122").Value = 4275.1665
$ws.Range("J122").Value = 7499.5
$ws.Range("L122").Value = 22498.5
$ws.Range("N122").Value = -27398.5
$ws.Range("H126").Value = 6229.5
$ws.Range("I126").Value = 6856.4287
$ws.Range("J126").Value = 4766.6665
$ws.Range("K126").Value = 20569.2861
$ws.Range("L126").Value = 14299.9995
$ws.Range("M126").Value = -18099.2861
$ws.Range("N126").Value = -19239.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3325.3333
$ws.Range("I16").Value = 3650.3333
$ws.Range("K16").Value = 3650.3333
$ws.Range("M16").Value = -3480.3333
$ws.Range("H22").Value = 854.3333
$ws.Range("I22").Value = 924.3333
$ws.Range("J22").Value = 761
$ws.Range("K22").Value = 924.3333
$ws.Range("L22").Value = 761
$ws.Range("M22").Value = -629.3333
$ws.Range("N22").Value = -1351
$ws.Range("H27").Value = 854.3333
$ws.Range("I27").Value = 924.3333
$ws.Range("J27").Value = 761
$ws.Range("K27").Value = 924.3333
$ws.Range("L27").Value = 761
$ws.Range("M27").Value = -817.3333
$ws.Range("N27").Value = -975
$ws.Range("H82").Value = 1356.2
$ws.Range("J82").Value = 1599.1666
$ws.Range("L82").Value = 1599.1666
$ws.Range("N82").Value = -2321.1666
$ws.Range("H85").Value = 1356.2
$ws.Range("J85").Value = 1599.1666
$ws.Range("L85").Value = 1599.1666
$ws.Range("N85").Value = -4095.1666
$ws.Range("H100").Value = 1801.5
$ws.Range("I100").Value = 1801.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1801.5
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1260.5
$ws.Range("N100").ClearContents()
$ws.Range("H106").Value = 7933.75
$ws.Range("J106").Value = 7933.75
$ws.Range("L106").Value = 7933.75
$ws.Range("N106").Value = -10457.75
$ws.Range("H122").Value = 3628.5
$ws.Range("I122").Value = 3503.2856
$ws.Range("J122").Value = 4505
$ws.Range("K122").Value = 10509.8568
$ws.Range("L122").Value = 13515
$ws.Range("M122").Value = -8059.856800000001
$ws.Range("N122").Value = -18415

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 733.3333
$ws.Range("I81").Value = 750
$ws.Range("J81").Value = 700
$ws.Range("K81").Value = 1500
$ws.Range("L81").Value = 1400
$ws.Range("M81").Value = -439
$ws.Range("N81").Value = -3522
$ws.Range("H84").Value = 733.3333
$ws.Range("I84").Value = 750
$ws.Range("J84").Value = 700
$ws.Range("K84").Value = 7500
$ws.Range("L84").Value = 7000
$ws.Range("M84").Value = -2196
$ws.Range("N84").Value = -17608
$ws.Range("H101").Value = 30957.143
$ws.Range("J101").Value = 30957.143
$ws.Range("L101").Value = 30957.143
$ws.Range("N101").Value = -37447.143
$ws.Range("H107").Value = 575.4
$ws.Range("I107").Value = 294.5
$ws.Range("K107").Value = 883.5
$ws.Range("M107").Value = 1036.5
$ws.Range("H113").Value = 8423.308000000001
$ws.Range("I113").Value = 14904.286
$ws.Range("K113").Value = 44712.858
$ws.Range("M113").Value = -42542.858
